$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 / Row 9 swap: Cardano <-> OKB (text columns, safe to set directly) ---
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'

# --- Plain text/value updates (Price values that are naturally text, and all Volume% cells) ---
$ws.Range("D2").Value = '30.086.33'
$ws.Range("E2").Value = '  +4.27%  '
$ws.Range("D3").Value = '1.907.72'
$ws.Range("E3").Value = '  +5.40%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +1.34%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +3.07%  '
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("E9").Value = '  +7.81%  '
$ws.Range("E10").Value = '  +6.07%  '
$ws.Range("D11").Value = '1.909.80'
$ws.Range("E11").Value = '  +5.50%  '
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("E13").Value = '  +4.47%  '
$ws.Range("E14").Value = '  +8.39%  '
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("E16").Value = '  +4.30%  '
$ws.Range("D17").Value = '30.105.41'
$ws.Range("E17").Value = '  +4.45%  '
$ws.Range("E18").Value = '  +10.95%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +6.42%  '
$ws.Range("D21").Value = '2.157.81'
$ws.Range("E21").Value = '  +5.45%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +5.96%  '
$ws.Range("E24").Value = '  +7.71%  '
$ws.Range("E25").Value = '  +4.22%  '
$ws.Range("E26").Value = '  +3.81%  '
$ws.Range("E27").Value = '  +2.82%  '
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("E29").Value = '  +6.35%  '
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("E32").Value = '  +5.92%  '
$ws.Range("E33").Value = '  +5.98%  '
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("E35").Value = '  +4.95%  '
$ws.Range("E36").Value = '  +6.70%  '
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("E41").Value = '  +6.88%  '
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("E43").Value = '  +5.72%  '
$ws.Range("E44").Value = '  +6.40%  '
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  +5.99%  '
$ws.Range("E47").Value = '  +4.96%  '
$ws.Range("E48").Value = '  +4.05%  '
$ws.Range("E49").Value = '  +5.89%  '
$ws.Range("E50").Value = '  +5.17%  '
$ws.Range("E51").Value = '  +5.40%  '

# --- Price values that look numeric: force text format first so Excel keeps them as strings ---
$numericLookingCells = @("D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D18","D19","D22","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D48","D49","D50","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").Value = '251.75'
$ws.Range("D6").Value = '0.9999'
$ws.Range("D7").Value = '0.5106'
$ws.Range("D8").Value = '45.12'
$ws.Range("D9").Value = '0.2996'
$ws.Range("D10").Value = '0.06807'
$ws.Range("D12").Value = '17.32'
$ws.Range("D13").Value = '0.07347'
$ws.Range("D14").Value = '0.7023'
$ws.Range("D15").Value = '86.73'
$ws.Range("D16").Value = '4.893'
$ws.Range("D18").Value = '0.000008173'
$ws.Range("D19").Value = '0.9992'
$ws.Range("D22").Value = '0.9988'
$ws.Range("D25").Value = '9.295'
$ws.Range("D26").Value = '147.81'
$ws.Range("D27").Value = '135.45'
$ws.Range("D28").Value = '17.09'
$ws.Range("D29").Value = '2.003'
$ws.Range("D30").Value = '1.400'
$ws.Range("D31").Value = '4.265'
$ws.Range("D32").Value = '0.08824'
$ws.Range("D33").Value = '4.013'
$ws.Range("D34").Value = '0.05073'
$ws.Range("D35").Value = '1.143'
$ws.Range("D36").Value = '0.7175'
$ws.Range("D37").Value = '2.688'
$ws.Range("D38").Value = '2.818'
$ws.Range("D39").Value = '2.263'
$ws.Range("D40").Value = '0.9648'
$ws.Range("D41").Value = '0.01700'
$ws.Range("D42").Value = '6.172'
$ws.Range("D43").Value = '0.4313'
$ws.Range("D44").Value = '105.56'
$ws.Range("D46").Value = '7.615'
$ws.Range("D48").Value = '0.05743'
$ws.Range("D49").Value = '33.35'
$ws.Range("D50").Value = '8.504'
$ws.Range("D51").Value = '0.3814'
